$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.07
$ws.Range("D2").Value = 0.0065194514888494903

# Row 3
$ws.Range("B3").Value = 0.17
$ws.Range("C3").Value = 0.39
$ws.Range("D3").Value = 0.017166580144957491
$ws.Range("F3").Value = 8840178.3152753748

# Row 5
$ws.Range("D5").Value = 0.019334296374540139
$ws.Range("E5").Value = 0.062840615108861048
$ws.Range("F5").Value = 9146002.7733741403
$ws.Range("G5").Value = 0.036962696997926363

# Row 6
$ws.Range("D6").Value = 0.017536085010108709
$ws.Range("E6").Value = 0.013967596463339759
$ws.Range("F6").Value = 9143534.0149630606
$ws.Range("G6").Value = 0.0082157062748250461

# Row 7
$ws.Range("E7").Value = 0.1113969968370139
$ws.Range("F7").Value = 8927469.1720317546
$ws.Range("G7").Value = 0.06552344265619546

# Row 8
$ws.Range("B8").Value = 0.04
$ws.Range("C8").Value = 0.28000000000000003
$ws.Range("E8").Value = 0.24690013833237789
$ws.Range("F8").Value = 8397441.337237617
$ws.Range("G8").Value = 0.1452260609816807

# Row 9
$ws.Range("B9").Value = 0.06
$ws.Range("C9").Value = 0.28999999999999998
$ws.Range("D9").Value = 0.018029877187196681
$ws.Range("E9").Value = 0.38828426439893238
$ws.Range("F9").Value = 7653183.6090507964
$ws.Range("G9").Value = 0.22838786013118931

# Row 10
$ws.Range("B10").Value = 0.07
$ws.Range("C10").Value = 0.3
$ws.Range("D10").Value = 0.01386146741153913
$ws.Range("E10").Value = 0.47410653732625269
$ws.Range("F10").Value = 7089451.9449234158
$ws.Range("G10").Value = 0.27886831237358911

# Update font family for bold header font
$ws.Range("A1:G1").Font.Name = "Calibri"

# Update selection
$ws.Range("F6").Select()
